$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view change ---
$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
$wb.Windows.Item(1).Left = 13640
$wb.Windows.Item(1).Top = -17960

# --- Insert a new row at 11 (shifts 11..33 down to 12..34), this is the
#     "temp" row that the rework introduces between D and Q ---
$ws.Rows("11:11").Insert()

# Row 8: D0 changes from 30e9 to 30e6
$ws.Range("D8").Value = 30000000

# Row 10: now holds the new "temp" label/value (previously the D formula row)
$ws.Range("C10").Value = "temp"
$ws.Range("D10").Value = 1000
$ws.Range("D10").Style = "Normal"

# Row 11 (newly inserted): now holds the D formula, referencing the new temp cell
$ws.Range("C11").Value = "D"
$ws.Range("D11").Formula = "=D8*EXP(-D9/8.31446/D10)"
$ws.Range("D11").NumberFormat = "0.0000E+00"

# Row 16 (was row 15, "dt"): value changes from 100000 to 10000000
$ws.Range("D16").Value = 10000000

# Row 19 (was row 18): first data row, formulas now reference the shifted cells
$ws.Range("C19").Formula = "=D16"
$ws.Range("D19").Formula = "=`$D`$12*(1-ERF(0,`$D`$14/SQRT(4*`$D`$11*C19)))"
$ws.Range("F19").Formula = "=ABS(E19-D19)/D19*100"

# Row 20 (was row 19): second data row, formulas now reference the shifted cells
$ws.Range("C20").Formula = "=C19+`$D`$16"
$ws.Range("D20").Formula = "=`$D`$12*(1-ERF(0,`$D`$14/SQRT(4*`$D`$11*C20)))"
$ws.Range("F20").Formula = "=ABS(E20-D20)/D20*100"

$wb.Application.Calculate()

# Selection ends up on the new D formula cell
$ws.Range("D11").Select()
